# New Password update SECR web & mobile
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update passwords (column B) for the SECR / ASM login rows.
# Order matters for shared-string allocation: ASM-005's password first,
# then SECR-008 (mobile) password, then SECR-008 (web) password.
$ws.Range("B7").Value = "Fosroc@7"
$ws.Range("B5").Value = "Fosroc@4"
$ws.Range("B11").Value = "Fosroc@0"

# Update the last selected cell on the active sheet.
[void]$ws.Range("J13").Select()
